$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update numeric cell values (rows 2-22) ---
$ws.Range("E2").Value = -0.139989
$ws.Range("K2").Value = 0.13161900000000001
$ws.Range("W2").Value = 0.128076
$ws.Range("E3").Value = 0.050680000000000003
$ws.Range("K3").Value = 0.161324
$ws.Range("W3").Value = 0.123608
$ws.Range("E4").Value = 0.50140499999999999
$ws.Range("G4").Value = 0.34389500000000001
$ws.Range("H4").Value = 0.162221
$ws.Range("K4").Value = 0.128216
$ws.Range("W4").Value = 0.0082529999999999999
$ws.Range("E5").Value = 0.59726999999999997
$ws.Range("G5").Value = 0.51836199999999999
$ws.Range("K5").Value = 0.20801900000000001
$ws.Range("W5").Value = 0.0051180000000000002
$ws.Range("E6").Value = 0.68357599999999996
$ws.Range("G6").Value = 0.37781999999999999
$ws.Range("H6").Value = 0.12236900000000001
$ws.Range("K6").Value = 0.135077
$ws.Range("W6").Value = 0.015656
$ws.Range("E7").Value = 0.73653299999999999
$ws.Range("G7").Value = 0.350157
$ws.Range("H7").Value = 0.103294
$ws.Range("K7").Value = 0.133048
$ws.Range("W7").Value = 0.030137000000000001
$ws.Range("E8").Value = 0.75483199999999995
$ws.Range("G8").Value = 0.31860300000000003
$ws.Range("H8").Value = 0.087749999999999995
$ws.Range("K8").Value = 0.13125800000000001
$ws.Range("W8").Value = 0.043320999999999998
$ws.Range("E9").Value = 0.78091500000000003
$ws.Range("G9").Value = 0.38547100000000001
$ws.Range("K9").Value = 0.14732500000000001
$ws.Range("W9").Value = 0.054004000000000003
$ws.Range("E10").Value = 0.72594899999999996
$ws.Range("G10").Value = 0.303618
$ws.Range("K10").Value = 0.12481200000000001
$ws.Range("W10").Value = 0.072359000000000007
$ws.Range("E11").Value = 0.72818799999999995
$ws.Range("G11").Value = 0.30460300000000001
$ws.Range("K11").Value = 0.126022
$ws.Range("W11").Value = 0.072172
$ws.Range("E12").Value = 0.70253500000000002
$ws.Range("G12").Value = 0.26140099999999999
$ws.Range("K12").Value = 0.15115000000000001
$ws.Range("W12").Value = 0.076502000000000001
$ws.Range("E13").Value = 0.70993499999999998
$ws.Range("G13").Value = 0.27426699999999998
$ws.Range("K13").Value = 0.152584
$ws.Range("W13").Value = 0.076065999999999995
$ws.Range("E14").Value = 0.68874999999999997
$ws.Range("G14").Value = 0.28818100000000002
$ws.Range("K14").Value = 0.13259799999999999
$ws.Range("W14").Value = 0.079020000000000007
$ws.Range("E15").Value = 0.66695000000000004
$ws.Range("G15").Value = 0.26180199999999998
$ws.Range("K15").Value = 0.12664600000000001
$ws.Range("W15").Value = 0.086053000000000004
$ws.Range("E16").Value = 0.624641
$ws.Range("G16").Value = 0.24431
$ws.Range("K16").Value = 0.10219200000000001
$ws.Range("W16").Value = 0.096202999999999997
$ws.Range("E17").Value = 0.580766
$ws.Range("G17").Value = 0.21432000000000001
$ws.Range("K17").Value = 0.104239
$ws.Range("W17").Value = 0.101949
$ws.Range("E18").Value = 0.50804099999999996
$ws.Range("G18").Value = 0.18524099999999999
$ws.Range("K18").Value = 0.098392999999999994
$ws.Range("W18").Value = 0.10946699999999999
$ws.Range("E19").Value = 0.41377399999999998
$ws.Range("G19").Value = 0.15631600000000001
$ws.Range("K19").Value = 0.090805999999999998
$ws.Range("W19").Value = 0.11706
$ws.Range("E20").Value = 0.48672599999999999
$ws.Range("G20").Value = 0.18254400000000001
$ws.Range("W20").Value = 0.135047
$ws.Range("E21").Value = 0.49192200000000003
$ws.Range("G21").Value = 0.190412
$ws.Range("W21").Value = 0.13528299999999999
$ws.Range("E22").Value = 0.49446299999999999
$ws.Range("G22").Value = 0.193795
$ws.Range("W22").Value = 0.13387299999999999

# --- Row 6 loses its custom "Neutral" style; row 12 gains it ---
$ws.Rows(6).ClearFormats()
foreach ($col in @("A","B","C","D","E","F","G","K","W")) {
    $ws.Range($col + "12").Style = "Neutral"
}

# --- Update the active selection on the sheet view ---
$ws.Range("I38").Select()
